$d = $word.ActiveDocument

# --- 1. Sol Da Silva replaces the placeholder game name in the title ---
$word.Application.UserName = "Sol Da Silva"
$d.TrackRevisions = $true

$d.Content.Find.Execute("<game name>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Puzzle Knight", 2)

# --- 2. Move the _GoBack bookmark so it sits right after the new title text ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$titlePara = $d.Paragraphs(3)
$bookmarkRange = $d.Range($titlePara.Range.End - 1, $titlePara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
